$d = $word.ActiveDocument

# The "Hoja De Datos" table has a header row, two filled data rows
# (Nery Javier..., Hector Mauricio...) and then several empty rows.
# Row 4 is the first fully empty row - fill it with Andre Figueroa's info.
$table = $d.Tables(1)
$row = $table.Rows(4)

$row.Cells(1).Range.Text = "André Sebastián Figueroa Barrios"
$row.Cells(2).Range.Text = "Desarrollador"
